$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.072064399719238
$ws.Range("B1").Value = 1.716970443725586
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.853051543235779
$ws.Range("E1").Value = 1.158736109733582
